$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 2555.25
$ws.Range("I20").Value = 2555.25
$ws.Range("K20").Value = 2555.25
$ws.Range("M20").Value = -2325.25

$ws.Range("H35").Value = 2555.25
$ws.Range("I35").Value = 2555.25
$ws.Range("K35").Value = 2555.25
$ws.Range("M35").Value = -2176.25

$ws.Range("H100").Value = 907.6087
$ws.Range("I100").Value = 1011.5625
$ws.Range("J100").Value = 670
$ws.Range("K100").Value = 1011.5625
$ws.Range("L100").Value = 670
$ws.Range("M100").Value = -470.5625
$ws.Range("N100").Value = -1752

$ws.Range("H129").Value = 905.7347
$ws.Range("J129").Value = 991.85
$ws.Range("L129").Value = 2975.55
$ws.Range("N129").Value = -12975.55

$ws.Range("H137").Value = 4350026.5
$ws.Range("I137").Value = 5557549
$ws.Range("J137").Value = 2946.6
$ws.Range("K137").Value = 16672647
$ws.Range("L137").Value = 8839.8
$ws.Range("M137").Value = -16670097
$ws.Range("N137").Value = -13939.8

$ws.Range("H138").Value = 3089286
$ws.Range("J138").Value = 4389563.5
$ws.Range("L138").Value = 13168690.5
$ws.Range("N138").Value = -13178970.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15905.482
$ws.Range("I32").Value = 16148.8
$ws.Range("J32").Value = 14080.6
$ws.Range("K32").Value = 16148.8
$ws.Range("L32").Value = 14080.6
$ws.Range("M32").Value = -15861.8
$ws.Range("N32").Value = -14654.6

$ws.Range("H45").Value = 2600
$ws.Range("I45").Value = 3966.6667
$ws.Range("J45").Value = 1780
$ws.Range("K45").Value = 3966.6667
$ws.Range("L45").Value = 1780
$ws.Range("M45").Value = -3589.6667
$ws.Range("N45").Value = -2534

$ws.Range("H61").Value = 35786900
$ws.Range("I61").Value = 45500790
$ws.Range("J61").Value = 169283.33
$ws.Range("K61").Value = 45500790
$ws.Range("L61").Value = 169283.33
$ws.Range("M61").Value = -45500578
$ws.Range("N61").Value = -169707.33

$ws.Range("H74").Value = 7877299
$ws.Range("I74").Value = 10460153
$ws.Range("J74").Value = 128737.5
$ws.Range("K74").Value = 10460153
$ws.Range("L74").Value = 128737.5
$ws.Range("M74").Value = -10459279
$ws.Range("N74").Value = -130485.5

$ws.Range("H77").Value = 7877299
$ws.Range("I77").Value = 10460153
$ws.Range("J77").Value = 128737.5
$ws.Range("K77").Value = 52300765
$ws.Range("L77").Value = 643687.5
$ws.Range("M77").Value = -52296397
$ws.Range("N77").Value = -652423.5

$ws.Range("H132").Value = 34462.95
$ws.Range("I132").Value = 23471.955
$ws.Range("K132").Value = 70415.865
$ws.Range("M132").Value = -67885.865

$ws.Range("H136").Value = 35786900
$ws.Range("I136").Value = 45500790
$ws.Range("J136").Value = 169283.33
$ws.Range("K136").Value = 136502370
$ws.Range("L136").Value = 507849.99
$ws.Range("M136").Value = -136499820
$ws.Range("N136").Value = -512949.99

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2420.3333
$ws.Range("I107").Value = 2660.1667
$ws.Range("J107").Value = 1940.6666
$ws.Range("K107").Value = 2660.1667
$ws.Range("L107").Value = 1940.6666
$ws.Range("M107").Value = -740.1667000000002
$ws.Range("N107").Value = -5780.6666

$ws.Range("H134").Value = 2558.8076
$ws.Range("I134").Value = 2002.025
$ws.Range("J134").Value = 4414.75
$ws.Range("K134").Value = 6006.075000000001
$ws.Range("L134").Value = 13244.25
$ws.Range("M134").Value = -3471.075000000001
$ws.Range("N134").Value = -18314.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 800
$ws.Range("J19").Value = 1250
$ws.Range("L19").Value = 1250
$ws.Range("N19").Value = -1590

$ws.Range("H24").Value = 800
$ws.Range("J24").Value = 1250
$ws.Range("L24").Value = 1250
$ws.Range("N24").Value = -1590

$ws.Range("H31").Value = 3149.5806
$ws.Range("I31").Value = 1602.1818
$ws.Range("J31").Value = 6932.1113
$ws.Range("K31").Value = 1602.1818
$ws.Range("L31").Value = 6932.1113
$ws.Range("M31").Value = -1307.1818
$ws.Range("N31").Value = -7522.1113

$ws.Range("H34").Value = 3149.5806
$ws.Range("I34").Value = 1602.1818
$ws.Range("J34").Value = 6932.1113
$ws.Range("K34").Value = 1602.1818
$ws.Range("L34").Value = 6932.1113
$ws.Range("M34").Value = -1400.1818
$ws.Range("N34").Value = -7336.1113

$ws.Range("H134").Value = 34540.914
$ws.Range("I134").Value = 2497.3845
$ws.Range("J134").Value = 127111.11
$ws.Range("K134").Value = 7492.1535
$ws.Range("L134").Value = 381333.33
$ws.Range("M134").Value = -4957.1535
$ws.Range("N134").Value = -386403.33

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 3900
$ws.Range("J102").Value = 3900
$ws.Range("L102").Value = 11700
$ws.Range("N102").Value = -16568

$ws.Range("H106").Value = 3489
$ws.Range("J106").Value = 3489
$ws.Range("L106").Value = 10467
$ws.Range("N106").Value = -12359

$ws.Range("H122").Value = 948.5217
$ws.Range("I122").Value = 359.5
$ws.Range("J122").Value = 1072.5264
$ws.Range("K122").Value = 3235.5
$ws.Range("L122").Value = 9652.7376
$ws.Range("M122").Value = -785.5
$ws.Range("N122").Value = -14552.7376

$ws.Range("H131").Value = 972.5161
$ws.Range("I131").Value = 380.8889
$ws.Range("J131").Value = 1072.9811
$ws.Range("K131").Value = 1142.6667
$ws.Range("L131").Value = 3218.9433
$ws.Range("M131").Value = 3897.3333
$ws.Range("N131").Value = -13298.9433

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 800
$ws.Range("I102").Value = 800
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 800
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 822
$ws.Range("N102").ClearContents()

$ws.Range("H132").Value = 46021.465
$ws.Range("I132").Value = 32544.219
$ws.Range("J132").Value = 79196.234
$ws.Range("K132").Value = 97632.657
$ws.Range("L132").Value = 237588.702
$ws.Range("M132").Value = -95102.657
$ws.Range("N132").Value = -242648.702

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 52891.668
$ws.Range("I132").Value = 24432.695
$ws.Range("J132").Value = 146399.72
$ws.Range("K132").Value = 73298.08499999999
$ws.Range("L132").Value = 439199.16
$ws.Range("M132").Value = -70768.08499999999
$ws.Range("N132").Value = -444259.16

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 68307.2
$ws.Range("I132").Value = 56684.055
$ws.Range("J132").Value = 85741.914
$ws.Range("K132").Value = 170052.165
$ws.Range("L132").Value = 257225.742
$ws.Range("M132").Value = -167522.165
$ws.Range("N132").Value = -262285.742

